$wb = $excel.ActiveWorkbook

# Sheet "About": A11 loses its bold styling (was style index 2 / bold font)
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A11").Font.Bold = $false

# Sheet "QSfHO": fix the quantization size for health outcomes data issue (was 0, now 1)
$wsQSfHO = $wb.Worksheets.Item("QSfHO")
$wsQSfHO.Range("B2").Value = 1

# Leave the cursor/selection on B3 of the QSfHO sheet (as it was when the
# workbook was last saved), but keep "About" as the active/visible sheet.
$wsQSfHO.Range("B3").Select() | Out-Null
$wsAbout.Activate() | Out-Null
